$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph, and pull in the blank
# paragraph right before it plus the copyright paragraph right after it.
# Together these three paragraphs form the site-footer block that the
# Jekyll rebuild dropped from the page.
$toDelete = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "Ver no Jupiter Salvar em pdf Salvar em docx") {
        $toDelete += $p.Previous()
        $toDelete += $p
        $toDelete += $p.Next()
    }
}

# Delete from the end backwards so earlier ranges in the list stay valid.
for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
    $toDelete[$i].Range.Delete()
}
